$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 104.794801
$ws.Range("H2").Value = 314.384403
$ws.Range("I2").Value = 0.3872421191355361
$ws.Range("J2").Value = 0.3872421191355361
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.127396333333333
$ws.Range("N2").Value = 6.382189
$ws.Range("O2").Value = 0.06137654768277986
$ws.Range("P2").Value = 0.06137654768277986
$ws.Range("Q2").Value = 222.9400753997964
$ws.Range("R2").Value = 2006.460678598167
$ws.Range("S2").Value = 0.02376758438990295
$ws.Range("T2").Value = 0.02376758438990295
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 104.794801
$ws.Range("H3").Value = 314.384403
$ws.Range("I3").Value = 0.3872421191355361
$ws.Range("J3").Value = 0.3872421191355361
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.721182333333333
$ws.Range("N3").Value = 11.163547
$ws.Range("O3").Value = 0.1073581454191429
$ws.Range("P3").Value = 0.1073581454191429
$ws.Range("Q3").Value = 389.9605621063823
$ws.Range("R3").Value = 3509.645058957441
$ws.Range("S3").Value = 0.04157359573856994
$ws.Range("T3").Value = 0.04157359573856993
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 104.794801
$ws.Range("H4").Value = 314.384403
$ws.Range("I4").Value = 0.3872421191355361
$ws.Range("J4").Value = 0.3872421191355361
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 28.81280933333333
$ws.Range("N4").Value = 86.43842799999999
$ws.Range("O4").Value = 0.8312653068980773
$ws.Range("P4").Value = 0.8312653068980772
$ws.Range("Q4").Value = 3019.432620337609
$ws.Range("R4").Value = 27174.89358303848
$ws.Range("S4").Value = 0.3219009390070632
$ws.Range("T4").Value = 0.3219009390070632
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 136.674446
$ws.Range("H5").Value = 410.023338
$ws.Range("I5").Value = 0.5050451128841343
$ws.Range("J5").Value = 0.5050451128841343
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.127396333333333
$ws.Range("N5").Value = 6.382189
$ws.Range("O5").Value = 0.06137654768277986
$ws.Range("P5").Value = 0.06137654768277986
$ws.Range("Q5").Value = 290.7607152807647
$ws.Range("R5").Value = 2616.846437526882
$ws.Range("S5").Value = 0.03099792545288801
$ws.Range("T5").Value = 0.03099792545288801
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 136.674446
$ws.Range("H6").Value = 410.023338
$ws.Range("I6").Value = 0.5050451128841343
$ws.Range("J6").Value = 0.5050451128841343
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.721182333333333
$ws.Range("N6").Value = 11.163547
$ws.Range("O6").Value = 0.1073581454191429
$ws.Range("P6").Value = 0.1073581454191429
$ws.Range("Q6").Value = 508.5905338733206
$ws.Range("R6").Value = 4577.314804859885
$ws.Range("S6").Value = 0.05422070667224232
$ws.Range("T6").Value = 0.05422070667224232
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 136.674446
$ws.Range("H7").Value = 410.023338
$ws.Range("I7").Value = 0.5050451128841343
$ws.Range("J7").Value = 0.5050451128841343
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 28.81280933333333
$ws.Range("N7").Value = 86.43842799999999
$ws.Range("O7").Value = 0.8312653068980773
$ws.Range("P7").Value = 0.8312653068980772
$ws.Range("Q7").Value = 3937.974753336962
$ws.Range("R7").Value = 35441.77278003265
$ws.Range("S7").Value = 0.419826480759004
$ws.Range("T7").Value = 0.419826480759004
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 29.14904533333333
$ws.Range("H8").Value = 87.447136
$ws.Range("I8").Value = 0.1077127679803296
$ws.Range("J8").Value = 0.1077127679803296
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.127396333333333
$ws.Range("N8").Value = 6.382189
$ws.Range("O8").Value = 0.06137654768277986
$ws.Range("P8").Value = 0.06137654768277986
$ws.Range("Q8").Value = 62.01157216230045
$ws.Range("R8").Value = 558.104149460704
$ws.Range("S8").Value = 0.006611037839988902
$ws.Range("T8").Value = 0.006611037839988902
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 29.14904533333333
$ws.Range("H9").Value = 87.447136
$ws.Range("I9").Value = 0.1077127679803296
$ws.Range("J9").Value = 0.1077127679803296
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.721182333333333
$ws.Range("N9").Value = 11.163547
$ws.Range("O9").Value = 0.1073581454191429
$ws.Range("P9").Value = 0.1073581454191429
$ws.Range("Q9").Value = 108.4689125279324
$ws.Range("R9").Value = 976.220212751392
$ws.Range("S9").Value = 0.01156384300833062
$ws.Range("T9").Value = 0.01156384300833062
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 29.14904533333333
$ws.Range("H10").Value = 87.447136
$ws.Range("I10").Value = 0.1077127679803296
$ws.Range("J10").Value = 0.1077127679803296
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 28.81280933333333
$ws.Range("N10").Value = 86.43842799999999
$ws.Range("O10").Value = 0.8312653068980773
$ws.Range("P10").Value = 0.8312653068980772
$ws.Range("Q10").Value = 839.865885438023
$ws.Range("R10").Value = 7558.792968942207
$ws.Range("S10").Value = 0.08953788713201005
$ws.Range("T10").Value = 0.08953788713201004
